# "All Expenses" sheet update
#  - Remove the duplicate/obsolete row 8 ("For finally.", Entertainment, 89.42, 2023-02-24),
#    which shifts every subsequent row up by one.
#  - Append two new expense rows at the bottom of the table:
#      33: March water bill    | 34.65 | Utilities | 2023-04-01
#      34: February Water Bill | 31.9  | Utilities | 2023-03-01

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 8 entirely; rows 9..33 shift up to become rows 8..32.
$ws.Rows.Item(8).Delete()

# Add the two new rows at the end of the table (now rows 33 and 34).
$ws.Cells.Item(33, 1).Value = "March water bill"
$ws.Cells.Item(33, 2).Value = 34.65
$ws.Cells.Item(33, 3).Value = "Utilities"
$ws.Cells.Item(33, 4).Value = 45017
$ws.Cells.Item(33, 4).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(34, 1).Value = "February Water Bill"
$ws.Cells.Item(34, 2).Value = 31.9
$ws.Cells.Item(34, 3).Value = "Utilities"
$ws.Cells.Item(34, 4).Value = 44986
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD"
